$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Status and Date values -----------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# Status: active -> draft
$wsMeta.Range("B6").Value = "draft"

# Date: 2023-05-12T12:33:13+00:00 -> 2023-08-01T16:12:28+00:00
$wsMeta.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# --- Re-assert the wrap/top alignment so it is flagged as "applied" ----
# (mirrors applyAlignment="true" being added to both cellXfs records that
# already carry <alignment vertical="top" wrapText="true"/>)
$wsOwl = $wb.Worksheets.Item("Include from hp.owl")

$headerRanges = @($wsMeta.Range("A1:B1"), $wsOwl.Range("A1:A1"))
$bodyRanges = @($wsMeta.Range("A2:B14"), $wsOwl.Range("A2:B4"))

foreach ($rng in $headerRanges) {
    $rng.VerticalAlignment = -4160
    $rng.WrapText = $true
}

foreach ($rng in $bodyRanges) {
    $rng.VerticalAlignment = -4160
    $rng.WrapText = $true
}
